# Apply the diff: O287 -> 2, fill R289/R290 with 0, and append rows 291-298.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) O287: 0 -> 2
$ws.Range("O287").Value = 2

# 2) R289 / R290: empty inline string -> numeric 0
$ws.Range("R289").Value = 0
$ws.Range("R290").Value = 0

# 3) Append new weekly rows 291-298.
$dateFmt = "YYYY-MM-DD HH:MM:SS"

$newRows = @(
    @{ Row=291; A=45474; B=1820;              C=1870.25;            D=1809.900024414062; E=1858.150024414062; F=1858.150024414062; G=2230801; H=2024; I=7; J=1;  K=0; L=0; M=0; N=27; O=0; P=0; Q=0 },
    @{ Row=292; A=45481; B=1855;              C=1947.300048828125;  D=1807.650024414062; E=1933.599975585938; F=1933.599975585938; G=2218210; H=2024; I=7; J=8;  K=0; L=0; M=0; N=28; O=0; P=0; Q=0 },
    @{ Row=293; A=45488; B=1945.050048828125; C=1960;               D=1803.050048828125; E=1818.699951171875; F=1818.699951171875; G=3426458; H=2024; I=7; J=15; K=0; L=0; M=0; N=29; O=1; P=0; Q=1 },
    @{ Row=294; A=45495; B=1802;              C=1842;               D=1751.099975585938; E=1812.199951171875; F=1812.199951171875; G=3607918; H=2024; I=7; J=22; K=0; L=0; M=0; N=30; O=0; P=0; Q=0 },
    @{ Row=295; A=45502; B=1822;              C=1869;               D=1785.349975585938; E=1798.699951171875; F=1798.699951171875; G=1724379; H=2024; I=7; J=29; K=0; L=0; M=0; N=31; O=0; P=0; Q=0 },
    @{ Row=296; A=45509; B=1782.75;           C=1794.199951171875;  D=1722.199951171875; E=1740.599975585938; F=1740.599975585938; G=2671120; H=2024; I=8; J=5;  K=0; L=0; M=0; N=32; O=0; P=0; Q=0 },
    @{ Row=297; A=45516; B=1753.099975585938; C=1754.75;            D=1705.449951171875; E=1750.800048828125; F=1750.800048828125; G=1426168; H=2024; I=8; J=12; K=0; L=0; M=0; N=33; O=0; P=0; Q=0 },
    @{ Row=298; A=45523; B=1757.949951171875; C=1809;                D=1751.699951171875; E=1804.050048828125; F=1804.050048828125; G=2722947; H=2024; I=8; J=19; K=0; L=0; M=0; N=34; O=0; P=0; Q=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A
    $ws.Range("A$row").NumberFormat = $dateFmt

    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    # R column stays blank (matches the untouched "backup" column for these rows).
}
